$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.8.1 -> 1.8.2
$ws.Range("B3").Value = "1.8.2"

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Experimental: true -> (cleared)
$ws.Range("B7").ClearContents()

# Date: 2023-10-31 -> 2025-11-18
# Entered as plain text (not auto-converted to a date serial): force a text
# number format, assign the value, then restore the original General-format
# style (copied from a neighboring cell) so the cell keeps its original
# style index.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2025-11-18"
$ws.Range("B5").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
